$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 13899.635
$ws.Range("J17").Value = 13899.635
$ws.Range("L17").Value = 41698.905
$ws.Range("N17").Value = -42034.905
$ws.Range("H40").Value = 71431340
$ws.Range("I40").Value = 2491
$ws.Range("J40").Value = 125002990
$ws.Range("K40").Value = 2491
$ws.Range("L40").Value = 125002990
$ws.Range("M40").Value = -2316
$ws.Range("N40").Value = -125003340
$ws.Range("H98").Value = 3473376
$ws.Range("I98").Value = 3907255
$ws.Range("K98").Value = 3907255
$ws.Range("M98").Value = -3905757
$ws.Range("H110").Value = 95860.336
$ws.Range("J110").Value = 95860.336
$ws.Range("L110").Value = 95860.336
$ws.Range("N110").Value = -104040.336
$ws.Range("H122").Value = 3473376
$ws.Range("I122").Value = 3907255
$ws.Range("K122").Value = 11721765
$ws.Range("M122").Value = -11719315
$ws.Range("H132").Value = 5729.875
$ws.Range("I132").Value = 2617.2856
$ws.Range("K132").Value = 7851.8568
$ws.Range("M132").Value = -5321.8568
$ws.Range("H137").Value = 1451.4286
$ws.Range("I137").Value = 932.2
$ws.Range("J137").Value = 2749.5
$ws.Range("K137").Value = 2796.6
$ws.Range("L137").Value = 8248.5
$ws.Range("M137").Value = -246.6000000000004
$ws.Range("N137").Value = -13348.5
$ws.Range("H138").Value = 3487.25
$ws.Range("I138").Value = 1640.5
$ws.Range("J138").Value = 5744.3887
$ws.Range("K138").Value = 4921.5
$ws.Range("L138").Value = 17233.1661
$ws.Range("M138").Value = 218.5
$ws.Range("N138").Value = -27513.1661
$ws.Range("H64").Value = 3999
$ws.Range("I64").Value = 3999
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3999
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("M64").Value = -3751
$ws.Range("H67").Value = 3999
$ws.Range("I67").Value = 3999
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3999
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("M67").Value = -3141
$ws.Range("H88").Value = 1003
$ws.Range("I88").Value = 1003
$ws.Range("K88").Value = 1003
$ws.Range("M88").Value = -597
$ws.Range("H91").Value = 1003
$ws.Range("I91").Value = 1003
$ws.Range("K91").Value = 1003
$ws.Range("M91").Value = 401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4093.1
$ws.Range("I32").Value = 4332.4375
$ws.Range("K32").Value = 4332.4375
$ws.Range("M32").Value = -4045.4375
$ws.Range("H61").Value = 8755575
$ws.Range("I61").Value = 11114384
$ws.Range("K61").Value = 11114384
$ws.Range("M61").Value = -11114172
$ws.Range("H92").Value = 86484.336
$ws.Range("J92").Value = 86484.336
$ws.Range("L92").Value = 86484.336
$ws.Range("N92").Value = -91476.336
$ws.Range("H132").Value = 1642103.9
$ws.Range("I132").Value = 2604.1052
$ws.Range("K132").Value = 7812.3156
$ws.Range("M132").Value = -5282.3156
$ws.Range("H136").Value = 8755575
$ws.Range("I136").Value = 11114384
$ws.Range("K136").Value = 33343152
$ws.Range("M136").Value = -33340602
$ws.Range("H69").Value = 250459
$ws.Range("J69").Value = 250459
$ws.Range("L69").Value = 250459
$ws.Range("N69").Value = -251957
$ws.Range("H72").Value = 250459
$ws.Range("J72").Value = 250459
$ws.Range("L72").Value = 751377
$ws.Range("N72").Value = -758865

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26727.143
$ws.Range("I86").Value = 33100.24
$ws.Range("J86").Value = 7607.857
$ws.Range("K86").Value = 33100.24
$ws.Range("L86").Value = 7607.857
$ws.Range("M86").Value = -31977.24
$ws.Range("N86").Value = -9853.857
$ws.Range("H89").Value = 26727.143
$ws.Range("I89").Value = 33100.24
$ws.Range("J89").Value = 7607.857
$ws.Range("K89").Value = 165501.2
$ws.Range("L89").Value = 38039.285
$ws.Range("M89").Value = -159885.2
$ws.Range("N89").Value = -49271.285
$ws.Range("H99").Value = 2835.1428
$ws.Range("I99").Value = 2399.4546
$ws.Range("K99").Value = 2399.4546
$ws.Range("M99").Value = -901.4546
$ws.Range("H134").Value = 6669320
$ws.Range("I134").Value = 2618.5454
$ws.Range("K134").Value = 7855.6362
$ws.Range("M134").Value = -5320.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2877.2666
$ws.Range("I58").Value = 2704.6155
$ws.Range("J58").Value = 3999.5
$ws.Range("K58").Value = 2704.6155
$ws.Range("L58").Value = 3999.5
$ws.Range("M58").Value = -2501.6155
$ws.Range("N58").Value = -4405.5
$ws.Range("H132").Value = 2927.037
$ws.Range("I132").Value = 2846.7144
$ws.Range("K132").Value = 8540.143199999999
$ws.Range("M132").Value = -6010.143199999999
$ws.Range("H134").Value = 2037.04
$ws.Range("I134").Value = 1844.6086
$ws.Range("K134").Value = 5533.825800000001
$ws.Range("M134").Value = -2998.825800000001
$ws.Range("H136").Value = 2877.2666
$ws.Range("I136").Value = 2704.6155
$ws.Range("J136").Value = 3999.5
$ws.Range("K136").Value = 8113.8465
$ws.Range("L136").Value = 11998.5
$ws.Range("M136").Value = -5563.8465
$ws.Range("N136").Value = -17098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 71428720
$ws.Range("J12").Value = 153.90909
$ws.Range("L12").Value = 461.72727
$ws.Range("N12").Value = -807.72727
$ws.Range("H44").Value = 4846.1875
$ws.Range("I44").Value = 399.16666
$ws.Range("J44").Value = 7514.4
$ws.Range("K44").Value = 1197.49998
$ws.Range("L44").Value = 22543.2
$ws.Range("M44").Value = -799.4999800000001
$ws.Range("N44").Value = -23339.2
$ws.Range("H92").Value = 37037396
$ws.Range("I92").Value = 166666910
$ws.Range("K92").Value = 500000730
$ws.Range("M92").Value = -499999482
$ws.Range("H97").Value = 1576.7
$ws.Range("J97").Value = 1408.625
$ws.Range("L97").Value = 4225.875
$ws.Range("N97").Value = -5217.875
$ws.Range("H116").Value = 12313
$ws.Range("I116").Value = 3014.5
$ws.Range("J116").Value = 21611.5
$ws.Range("K116").Value = 9043.5
$ws.Range("L116").Value = 64834.5
$ws.Range("M116").Value = -5601.5
$ws.Range("N116").Value = -71718.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4389.4287
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 4620.1665
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 4620.1665
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -6616.1665
$ws.Range("H83").Value = 4389.4287
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 4620.1665
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 23100.8325
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -33084.8325
$ws.Range("H97").Value = 668.8684
$ws.Range("I97").Value = 582.931
$ws.Range("J97").Value = 945.7778
$ws.Range("K97").Value = 582.931
$ws.Range("L97").Value = 945.7778
$ws.Range("M97").Value = -86.93100000000004
$ws.Range("N97").Value = -1937.7778
$ws.Range("H102").Value = 2842.5454
$ws.Range("I102").Value = 2837
$ws.Range("K102").Value = 2837
$ws.Range("M102").Value = -1215
$ws.Range("H132").Value = 7694808
$ws.Range("I132").Value = 2708.75
$ws.Range("K132").Value = 8126.25
$ws.Range("M132").Value = -5596.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2785.5938
$ws.Range("I132").Value = 1820.0741
$ws.Range("K132").Value = 5460.2223
$ws.Range("M132").Value = -2930.2223
$ws.Range("H136").Value = 3408.7083
$ws.Range("I136").Value = 2400.8235
$ws.Range("J136").Value = 5856.4287
$ws.Range("K136").Value = 7202.470499999999
$ws.Range("L136").Value = 17569.2861
$ws.Range("M136").Value = -4652.470499999999
$ws.Range("N136").Value = -22669.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 99998.5
$ws.Range("J111").Value = 99998.5
$ws.Range("L111").Value = 99998.5
$ws.Range("N111").Value = -108178.5
$ws.Range("H122").Value = 2970.2083
$ws.Range("J122").Value = 3768.3635
$ws.Range("L122").Value = 11305.0905
$ws.Range("N122").Value = -16205.0905
$ws.Range("H126").Value = 6812.3447
$ws.Range("I126").Value = 7870.4346
$ws.Range("J126").Value = 2756.3333
$ws.Range("K126").Value = 23611.3038
$ws.Range("L126").Value = 8268.999899999999
$ws.Range("M126").Value = -21141.3038
$ws.Range("N126").Value = -13208.9999
$ws.Range("H136").Value = 280430.75
$ws.Range("J136").Value = 2004281.6
$ws.Range("L136").Value = 6012844.800000001
$ws.Range("N136").Value = -6017944.800000001
